$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "Student_ID"
$ws.Range("B1").Value = "First_Name"
$ws.Range("C1").Value = "Middle_Initial"
$ws.Range("D1").Value = "Last_Name"
$ws.Range("E1").Value = "Date_of_Birth"
$ws.Range("F1").Value = "CAIR_Program_Name"
$ws.Range("G1").Value = "Notes"

# Apply font formatting (Calibri 12) to A1:E1
$fmtRange = $ws.Range("A1:E1")
$fmtRange.Font.Name = "Calibri"
$fmtRange.Font.Size = 12

# Row height for header row
$ws.Rows.Item(1).RowHeight = 15.75

# Selection matches target sheetView
$ws.Range("P7").Select()
